$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Events_Calendar")

# --- Fix the Acquisition COI event description (row 4 / column E) ---
# "guest speaker" -> "guest speakers ... (ODASA(P) and C" (truncated import text)
$ws.Range("E4").Value = "Join the ACT-IAC Acquisition COI member meeting featuring guest speakers Elizabeth (Liz) Chirico who leads the Acquisition Innovation through Technology team in the Office of the Deputy Assistant Secretary of the Army (Procurement) (ODASA(P) and C"

# --- Realign the "Learn More" hyperlinks for rows 13-19 ---
# The links had drifted out of sync with their events (row 13's link actually
# belonged to row 14, row 14's to row 15, and row 15 had no link at all).
# Rebuild all hyperlinks on the sheet so everything lines back up:
#   row 13 "SBA on ISO, CMMI, and CMMC"       -> no real link yet (TBD)
#   row 14 "Climate Change Summit"            -> climate-change-summit
#   row 15 "ACT-IAC Cybersecurity COI Apr 23" -> act-iac-cybersecurity-coi-april-2023
#   rows 16-19 keep their (already-correct) links
$ws.Hyperlinks.Delete()

$ws.Range("F13").Value = "TBD"
$ws.Range("F13").Style = "Normal"

$ws.Hyperlinks.Add($ws.Range("F2"), "https://www.actiac.org/act-iac-event/federal-insights-exchange-session-featuring-dhs-pil")
$ws.Range("F2").Style = "Hyperlink"

$ws.Hyperlinks.Add($ws.Range("F3"), "https://www.actiac.org/act-iac-event/emerging-technology-and-innovation-opportunities-small-businesses")
$ws.Range("F3").Style = "Hyperlink"

$ws.Hyperlinks.Add($ws.Range("F4"), "https://www.actiac.org/act-iac-event/act-iac-acquisition-coi-march-2023")
$ws.Range("F4").Style = "Hyperlink"

$ws.Hyperlinks.Add($ws.Range("F5"), "https://www.actiac.org/act-iac-event/act-iac-it-management-and-modernization-coi-march-2023")
$ws.Range("F5").Style = "Hyperlink"

$ws.Hyperlinks.Add($ws.Range("F6"), "https://www.actiac.org/act-iac-event/voyagers-got-talent-2023")
$ws.Range("F6").Style = "Hyperlink"

$ws.Hyperlinks.Add($ws.Range("F7"), "https://www.actiac.org/act-iac-event/act-iac-cybersecurity-coi-march-2023")
$ws.Range("F7").Style = "Hyperlink"

$ws.Hyperlinks.Add($ws.Range("F8"), "https://www.actiac.org/node/8906556")
$ws.Range("F8").Style = "Hyperlink"

$ws.Hyperlinks.Add($ws.Range("F9"), "https://www.actiac.org/act-iac-event/federal-insights-exchange-session-feat-department-education-0")
$ws.Range("F9").Style = "Hyperlink"

$ws.Hyperlinks.Add($ws.Range("F10"), "https://www.actiac.org/act-iac-event/act-iac-evolving-workforce-coi-march-2023")
$ws.Range("F10").Style = "Hyperlink"

$ws.Hyperlinks.Add($ws.Range("F11"), "https://www.actiac.org/act-iac-event/act-iac-emerging-technology-coi-emerging-technology-et-accelerator-small-business")
$ws.Range("F11").Style = "Hyperlink"

$ws.Hyperlinks.Add($ws.Range("F12"), "https://www.actiac.org/act-iac-event/federal-insights-exchange-featuring-dept-transportation-cio-cordell-schachter")
$ws.Range("F12").Style = "Hyperlink"

$ws.Range("F14").Value = "https://www.actiac.org/act-iac-event/climate-change-summit"
$ws.Hyperlinks.Add($ws.Range("F14"), "https://www.actiac.org/act-iac-event/climate-change-summit")
$ws.Range("F14").Style = "Hyperlink"

$ws.Range("F15").Value = "https://www.actiac.org/act-iac-event/act-iac-cybersecurity-coi-april-2023"
$ws.Hyperlinks.Add($ws.Range("F15"), "https://www.actiac.org/act-iac-event/act-iac-cybersecurity-coi-april-2023")
$ws.Range("F15").Style = "Hyperlink"

$ws.Hyperlinks.Add($ws.Range("F16"), "https://www.actiac.org/act-iac-event/fie-session-featuring-gsa-ogp-office-evidence-and-analysis")
$ws.Range("F16").Style = "Hyperlink"

$ws.Hyperlinks.Add($ws.Range("F17"), "https://www.actiac.org/act-iac-event/emerging-technology-and-innovation-conference-2023")
$ws.Range("F17").Style = "Hyperlink"

$ws.Hyperlinks.Add($ws.Range("F18"), "https://www.actiac.org/act-iac-event/act-iac-cybersecurity-coi-may-2023")
$ws.Range("F18").Style = "Hyperlink"

$ws.Hyperlinks.Add($ws.Range("F19"), "https://www.actiac.org/act-iac-event/imagine-nation-elc-2023")
$ws.Range("F19").Style = "Hyperlink"

# --- Remove the leftover "Projects" sheet (no longer part of the workbook) ---
$wsProjects = $wb.Worksheets.Item("Projects")
$wsProjects.Delete()
